$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 4 per the FlashScore data refresh
$ws.Range("G4").Value = 2.9
$ws.Range("I4").Value = 2.25
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 2.12
$ws.Range("L4").Value = 2.85
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 3.35
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 1.91
$ws.Range("U4").Value = 1.57
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 10.75
$ws.Range("X4").Value = 16.5
$ws.Range("Y4").Value = 10.25
$ws.Range("AB4").Value = 27
$ws.Range("AC4").Value = 11.5
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 8.75
$ws.Range("AJ4").Value = 9
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 17.5
$ws.Range("AO4").Value = 15
$ws.Range("AQ4").Value = 65
$ws.Range("AR4").Value = 90
$ws.Range("AT4").Value = 2.82
$ws.Range("AW4").Value = 4.25
$ws.Range("AX4").Value = 11.75
$ws.Range("AY4").Value = 19
$ws.Range("BA4").Value = 75
